# Daily attendance processing - normalize the "Recorded By" (column G) values.
# For every data row, the list of recorders stored in column G (a comma
# separated list, e.g. "dnasr281@gmail.com, System") is re-sorted into
# ordinal (case-sensitive) alphabetical order, e.g. "System, dnasr281@gmail.com".

function Sort-Ordinal($items) {
    # simple insertion sort using ordinal (case-sensitive) string comparison,
    # since the built-in Sort-Object does not use ordinal comparison here.
    $list = @($items)
    $n = $list.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $list[$i]
        $j = $i - 1
        while ($j -ge 0 -and $list[$j].CompareTo($key) -gt 0) {
            $list[$j + 1] = $list[$j]
            $j = $j - 1
        }
        $list[$j + 1] = $key
    }
    return $list
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the 7th column ("Recorded By"); row 1 is the header row.
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = @($val -split "," | ForEach-Object { $_.Trim() })
    if ($parts.Count -le 1) { continue }

    $sorted = Sort-Ordinal $parts
    $newVal = [string]::Join(", ", $sorted)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}

Write-Host "Recorded By column normalized."
